# full_production_plan_report.xlsx
# fix(module3): use uncon_planned_qty for future production; keep produced for today
#
# Rewrites rows 2-9 with refreshed quantities/dates (the plan shifted one day
# forward and quantities were recomputed), and appends rows 10-13 so the
# report now covers the full simulated week.
#
# Columns: A produced_qty, B material, C location, D line,
#          E simulation_date, F production_plan_date, G available_date,
#          H uncon_planned_qty, I con_planned_qty, J changeover_id

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2;  A=675; B="MAT_A"; D="LINE_A"; Sim=45293; H=710; I=710 },
    @{ Row=3;  A=92;  B="MAT_B"; D="LINE_B"; Sim=45293; H=104; I=104 },
    @{ Row=4;  A=621; B="MAT_A"; D="LINE_A"; Sim=45294; H=660; I=660 },
    @{ Row=5;  A=75;  B="MAT_B"; D="LINE_B"; Sim=45294; H=80;  I=80  },
    @{ Row=6;  A=675; B="MAT_A"; D="LINE_A"; Sim=45295; H=710; I=710 },
    @{ Row=7;  A=92;  B="MAT_B"; D="LINE_B"; Sim=45295; H=104; I=104 },
    @{ Row=8;  A=817; B="MAT_A"; D="LINE_A"; Sim=45296; H=860; I=860 },
    @{ Row=9;  A=99;  B="MAT_B"; D="LINE_B"; Sim=45296; H=112; I=112 },
    @{ Row=10; A=798; B="MAT_A"; D="LINE_A"; Sim=45297; H=840; I=840 },
    @{ Row=11; A=99;  B="MAT_B"; D="LINE_B"; Sim=45297; H=112; I=112 },
    @{ Row=12; A=827; B="MAT_A"; D="LINE_A"; Sim=45298; H=870; I=870 },
    @{ Row=13; A=106; B="MAT_B"; D="LINE_B"; Sim=45298; H=120; I=120 }
)

foreach ($r in $rows) {
    $row = $r.Row

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = "PLANT_001"
    $ws.Cells.Item($row, 4).Value = $r.D

    # simulation_date / production_plan_date / available_date are 3
    # consecutive days, stored as date-formatted serials.
    $ws.Cells.Item($row, 5).Value = $r.Sim
    $ws.Cells.Item($row, 6).Value = $r.Sim + 1
    $ws.Cells.Item($row, 7).Value = $r.Sim + 2

    $ws.Cells.Item($row, 8).Value = $r.H
    $ws.Cells.Item($row, 9).Value = $r.I

    # Newly appended rows need the date number format + an (empty) changeover
    # id cell explicitly materialised, matching the existing rows' layout.
    if ($row -gt 9) {
        $ws.Cells.Item($row, 5).NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $ws.Cells.Item($row, 6).NumberFormat = "YYYY-MM-DD HH:MM:SS"
        $ws.Cells.Item($row, 7).NumberFormat = "YYYY-MM-DD HH:MM:SS"

        $ws.Cells.Item($row, 10).Value = "'"
        $ws.Cells.Item($row, 10).Style = "Normal"
    }
}
